$wb = $excel.ActiveWorkbook

# --- Locate the source sheet that the new "EXE Authorization" tab is based on ---
$src = $wb.Worksheets.Item("EXE Login")

# --- Duplicate "EXE Login" right after itself, producing the new third sheet. ---
# Worksheet.Copy faithfully reproduces cell values, styles, column widths,
# number formats and merged cells, which a plain value/format paste would not.
$src.Copy([System.Reflection.Missing]::Value, $src)

# The copy becomes the new active sheet, placed immediately after "EXE Login".
$newSheet = $wb.Worksheets.Item(3)
$newSheet.Name = "EXE Authorization"

# Keep the freshly-added sheet selected/active (matches activeTab moving to
# the new 3rd tab, 0-indexed = 2, and tabSelected="1" landing on the new sheet).
$newSheet.Activate()
$newSheet.Select()
